$wb = $excel.ActiveWorkbook

# RateSheetManagement sheet: rate sheet values change from "Schedule A" to "DRC - Original"
$wsRate = $wb.Worksheets.Item("RateSheetManagement")
$wsRate.Range("B2:B5").Value = "DRC - Original"
$wsRate.Range("B2:B5").Style = "Normal"

# StaffMember sheet: last row's Title cell formatting reset to default style
$wsStaff = $wb.Worksheets.Item("StaffMember")
$wsStaff.Range("A5").Style = "Normal"

# Make RateSheetManagement the active/selected sheet and cell, matching the
# updated workbook view state (tab focus moved away from SummaryLogs).
$wsRate.Activate()
$wsRate.Range("B2").Select()
